$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helper: find the paragraph whose trimmed text equals $text
# ------------------------------------------------------------------
function Find-Paragraph($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

# ------------------------------------------------------------------
# 1. Insert the new "Sound" section (heading + body paragraph) right
#    before the existing "Menu:" heading paragraph.
# ------------------------------------------------------------------
$menuPara = Find-Paragraph("Menu:")

# Insert two new (empty) paragraphs right before "Menu:". Both inherit
# Menu's paragraph formatting (ListParagraph style + numPr numbering),
# which is exactly what we want for the "Sound" heading paragraph.
$menuPara.Range.InsertParagraphBefore()
$menuPara.Range.InsertParagraphBefore()

# Re-locate "Menu:" (its paragraph index shifted down by two).
$menuPara = Find-Paragraph("Menu:")
$menuStart = $menuPara.Range.Start

# The two freshly inserted empty paragraphs sit right before Menu's
# start. The first one (closer to "Enemies" section) becomes the
# "Sound" heading.
$soundHeading = $d.Paragraphs.Item(1)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -lt $menuStart -and $p.Range.End -le $menuStart -and $p.Range.Text.TrimEnd([char]13, [char]7) -eq "") {
        $soundHeading = $p
        break
    }
}
$soundHeading.Range.Text = "Sound"

# Re-locate the (now second) empty paragraph that must hold the body
# text of the "Sound" section, and the "Menu:" heading again.
$menuPara = Find-Paragraph("Menu:")
$menuStart = $menuPara.Range.Start
$soundBody = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -lt $menuStart -and $p.Range.End -le $menuStart -and $p.Range.Text.TrimEnd([char]13, [char]7) -eq "") {
        $soundBody = $p
        break
    }
}

# Replace the body paragraph's contents (and paragraph formatting) with
# the required text, including the grammar-check markers around
# "shoots" that mirror the style used elsewhere in the document.
$xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">We are going to implement sound when the player ship </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>shoots</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> or the enemy ship is destroyed or when the boss comes to the game and also a menu background music. We also want to implement a little melody in the game sound.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$soundBody.Range.InsertXML($xml)

# ------------------------------------------------------------------
# 2. Resize the ASCII-art picture (uniform scale, aspect ratio locked)
#    by setting its width - height follows automatically.
# ------------------------------------------------------------------
$shape = $d.InlineShapes.Item(1)
$shape.Width = 255.6
